$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new header cells AD1:AF1 ("Wins", "Losses", "Ties") ---
# Copy the formatting of an existing header cell (AC1) so the new headers
# match the bold / bordered / centered style used by the rest of row 1.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Fill the team record (Wins/Losses/Ties) for every data row 2-60 ---
$firstRow = 2
$lastRow = 60
$rowCount = $lastRow - $firstRow + 1

$rng = $ws.Range("AD" + $firstRow + ":AF" + $lastRow)
$data = New-Object 'object[,]' $rowCount,3
for ($i = 0; $i -lt $rowCount; $i++) {
    $data[$i,0] = 56
    $data[$i,1] = 106
    $data[$i,2] = 0
}
$rng.Value = $data
